$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 351 (shifts old rows 351-383 down to 352-384)
$ws.Rows.Item(351).Insert()

# Populate the new row 351 with a copy of the (old) row 351 data, but with
# an updated Fecha (column D): 45033 -> 45212
$ws.Range("A351").Value = 5
$ws.Range("B351").Value = "Macroferia Regional de Talca"
$ws.Range("C351").Value = "Maule"
$ws.Range("D351").Value = 45212
$ws.Range("E351").Value = 7
$ws.Range("F351").Value = 100112017
$ws.Range("G351").Value = "Apio"
$ws.Range("H351").Value = "Americana (o)"
$ws.Range("I351").Value = "Primera"
$ws.Range("J351").Value = 500
$ws.Range("K351").Value = 7500
$ws.Range("L351").Value = 7500
$ws.Range("M351").Value = 7500
$ws.Range("N351").Value = "`$/docena de matas"
$ws.Range("O351").Value = "Provincia del Elquí"
$ws.Range("P351").Value = 1250
$ws.Range("Q351").Value = 6
$ws.Range("R351").Value = "Hortaliza"
